# Gauss-Seidel results table: the underlying system being solved changed
# (3 equations instead of 2, different starting point), so the iteration
# log is replaced:
#   - rows 2-4 get the new xn / error values
#   - rows 5-16 (the extra iterations the old 2-variable system needed to
#     converge) are removed entirely, since the corrected system converges
#     in 3 iterations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $cellAddr as a literal TEXT value (matching the
# original inline-string cells) without leaving any NumberFormat/style
# residue on the destination cell. A plain `Range.Value = "1.0"` would be
# silently reinterpreted as a number by Excel, and pre-formatting the
# destination as Text would permanently change its style (not part of this
# edit). Instead, stage the value as a text-producing formula in a scratch
# cell (row 16, which gets deleted at the end anyway), then copy/paste
# *values only* so the destination keeps its original default style.
$scratch = $ws.Range("Z16")
function Set-TextValue([string]$cellAddr, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)  # xlPasteValues
}

# Row 2 (iteration 1)
Set-TextValue "B2" "[1;3;6]"
Set-TextValue "C2" "1.0"

# Row 3 (iteration 2)
Set-TextValue "B3" "[4.39;2.895;6.82125]"
Set-TextValue "C3" "0.772209567198178"

# Row 4 (iteration 3)
Set-TextValue "B4" "[4.8591625;2.88320625;6.9355921875]"
Set-TextValue "C4" "0.0965521321832722"

$scratch.ClearContents()

# Drop the now-obsolete rows 5-16.
$ws.Rows("5:16").Delete()
